$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC4 is entered into the next empty row (A3), matching the prior TC2/TC3 entries.
$ws.Range("A3").Value = "TC4"

# Selection moves on to the next empty cell, as Excel leaves it after manual entry.
$ws.Range("A4").Select() | Out-Null
